$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Actualizar el texto de la tarea (shared string usado en B2)
$ws.Range("B2").Value = "Instalar y configurar Ruby On Rails."

# Cambiar la alineación vertical de "bottom" a "top" para las celdas
# combinadas D4:D6 y F4:F6 (estilo con borde que antes tenía vertical=bottom)
$xlVAlignTop = -4160
$ws.Range("D4:D6").VerticalAlignment = $xlVAlignTop
$ws.Range("F4:F6").VerticalAlignment = $xlVAlignTop

# Mover la selección activa a F10
[void]$ws.Range("F10").Select()
